$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 32.00264033333334
$ws.Range("H2").Value = 96.00792100000001
$ws.Range("I2").Value = 0.02419001798940439
$ws.Range("J2").Value = 0.02433952891158457
$ws.Range("M2").Value = 11.53121
$ws.Range("N2").Value = 34.59363
$ws.Range("O2").Value = 0.04670311854310271
$ws.Range("P2").Value = 0.04697417460197403
$ws.Range("Q2").Value = 369.0291662381367
$ws.Range("R2").Value = 3321.26249614323
$ws.Range("S2").Value = 0.00112974927771894
$ws.Range("T2").Value = 0.001143329280822569

$ws.Range("G3").Value = 32.00264033333334
$ws.Range("H3").Value = 96.00792100000001
$ws.Range("I3").Value = 0.02419001798940439
$ws.Range("J3").Value = 0.02433952891158457
$ws.Range("O3").Value = 0.04941252241252041
$ws.Range("P3").Value = 0.04969930333854504
$ws.Range("Q3").Value = 390.4377805260847
$ws.Range("R3").Value = 3513.940024734762
$ws.Range("S3").Value = 0.001195289806060716
$ws.Range("T3").Value = 0.001209657630494129

$ws.Range("G4").Value = 32.00264033333334
$ws.Range("H4").Value = 96.00792100000001
$ws.Range("I4").Value = 0.02419001798940439
$ws.Range("J4").Value = 0.02433952891158457
$ws.Range("M4").Value = 113.969907
$ws.Range("N4").Value = 341.909721
$ws.Range("O4").Value = 0.4615951038067463
$ws.Range("P4").Value = 0.4642741144067919
$ws.Range("Q4").Value = 3647.33794254445
$ws.Range("R4").Value = 32826.04148290004
$ws.Range("S4").Value = 0.01116599386490618
$ws.Range("T4").Value = 0.01130021323050443

$ws.Range("G5").Value = 32.00264033333334
$ws.Range("H5").Value = 96.00792100000001
$ws.Range("I5").Value = 0.02419001798940439
$ws.Range("J5").Value = 0.02433952891158457
$ws.Range("M5").Value = 4.2741545
$ws.Range("N5").Value = 8.548309
$ws.Range("O5").Value = 0.01731096253429049
$ws.Range("P5").Value = 0.01160762138918714
$ws.Range("Q5").Value = 136.7842291925982
$ws.Range("R5").Value = 820.7053751555891
$ws.Range("S5").Value = 0.0004187524951183923
$ws.Range("T5").Value = 0.000282524036396848

$ws.Range("G6").Value = 32.00264033333334
$ws.Range("H6").Value = 96.00792100000001
$ws.Range("I6").Value = 0.02419001798940439
$ws.Range("J6").Value = 0.02433952891158457
$ws.Range("M6").Value = 104.9290516666667
$ws.Range("N6").Value = 314.787155
$ws.Range("O6").Value = 0.4249782927033401
$ws.Range("P6").Value = 0.4274447862635018
$ws.Range("Q6").Value = 3358.006701006084
$ws.Range("R6").Value = 30222.06030905476
$ws.Range("S6").Value = 0.01028023254560016
$ws.Range("T6").Value = 0.01040380473336659

$ws.Range("I7").Value = 0.07580145430919519
$ws.Range("J7").Value = 0.07626995934880827
$ws.Range("M7").Value = 11.53121
$ws.Range("N7").Value = 34.59363
$ws.Range("O7").Value = 0.04670311854310271
$ws.Range("P7").Value = 0.04697417460197403
$ws.Range("Q7").Value = 1156.38390577523
$ws.Range("R7").Value = 10407.45515197707
$ws.Range("S7").Value = 0.003540164306341927
$ws.Range("T7").Value = 0.003582718387336382

$ws.Range("I8").Value = 0.07580145430919519
$ws.Range("J8").Value = 0.07626995934880827
$ws.Range("O8").Value = 0.04941252241252041
$ws.Range("P8").Value = 0.04969930333854504
$ws.Range("S8").Value = 0.003745541059954749
$ws.Range("T8").Value = 0.003790563845294922

$ws.Range("I9").Value = 0.07580145430919519
$ws.Range("J9").Value = 0.07626995934880827
$ws.Range("M9").Value = 113.969907
$ws.Range("N9").Value = 341.909721
$ws.Range("O9").Value = 0.4615951038067463
$ws.Range("P9").Value = 0.4642741144067919
$ws.Range("Q9").Value = 11429.23996679444
$ws.Range("R9").Value = 102863.15970115
$ws.Range("S9").Value = 0.03498958017055529
$ws.Range("T9").Value = 0.03541016783250998

$ws.Range("I10").Value = 0.07580145430919519
$ws.Range("J10").Value = 0.07626995934880827
$ws.Range("M10").Value = 4.2741545
$ws.Range("N10").Value = 8.548309
$ws.Range("O10").Value = 0.01731096253429049
$ws.Range("P10").Value = 0.01160762138918714
$ws.Range("Q10").Value = 428.6248775797835
$ws.Range("R10").Value = 2571.749265478701
$ws.Range("S10").Value = 0.001312196135591211
$ws.Range("T10").Value = 0.0008853128114896609

$ws.Range("I11").Value = 0.07580145430919519
$ws.Range("J11").Value = 0.07626995934880827
$ws.Range("M11").Value = 104.9290516666667
$ws.Range("N11").Value = 314.787155
$ws.Range("O11").Value = 0.4249782927033401
$ws.Range("P11").Value = 0.4274447862635018
$ws.Range("Q11").Value = 10522.59620591342
$ws.Range("R11").Value = 94703.36585322079
$ws.Range("S11").Value = 0.03221397263675201
$ws.Range("T11").Value = 0.03260119647217732

$ws.Range("G12").Value = 473.968811
$ws.Range("H12").Value = 1421.906433
$ws.Range("I12").Value = 0.3582615042098434
$ws.Range("J12").Value = 0.360475806319893
$ws.Range("M12").Value = 11.53121
$ws.Range("N12").Value = 34.59363
$ws.Range("O12").Value = 0.04670311854310271
$ws.Range("P12").Value = 0.04697417460197403
$ws.Range("Q12").Value = 5465.43389309131
$ws.Range("R12").Value = 49188.90503782179
$ws.Range("S12").Value = 0.01673192950054261
$ws.Range("T12").Value = 0.01693305346585803

$ws.Range("G13").Value = 473.968811
$ws.Range("H13").Value = 1421.906433
$ws.Range("I13").Value = 0.3582615042098434
$ws.Range("J13").Value = 0.360475806319893
$ws.Range("O13").Value = 0.04941252241252041
$ws.Range("P13").Value = 0.04969930333854504
$ws.Range("Q13").Value = 5782.501964773113
$ws.Range("R13").Value = 52042.51768295803
$ws.Range("S13").Value = 0.01770260460631216
$ws.Range("T13").Value = 0.01791539644449897

$ws.Range("G14").Value = 473.968811
$ws.Range("H14").Value = 1421.906433
$ws.Range("I14").Value = 0.3582615042098434
$ws.Range("J14").Value = 0.360475806319893
$ws.Range("M14").Value = 113.969907
$ws.Range("N14").Value = 341.909721
$ws.Range("O14").Value = 0.4615951038067463
$ws.Range("P14").Value = 0.4642741144067919
$ws.Range("Q14").Value = 54018.18131057057
$ws.Range("R14").Value = 486163.6317951352
$ws.Range("S14").Value = 0.1653717562257037
$ws.Range("T14").Value = 0.1673595857442426

$ws.Range("G15").Value = 473.968811
$ws.Range("H15").Value = 1421.906433
$ws.Range("I15").Value = 0.3582615042098434
$ws.Range("J15").Value = 0.360475806319893
$ws.Range("M15").Value = 4.2741545
$ws.Range("N15").Value = 8.548309
$ws.Range("O15").Value = 0.01731096253429049
$ws.Range("P15").Value = 0.01160762138918714
$ws.Range("Q15").Value = 2025.8159263953
$ws.Range("R15").Value = 12154.8955583718
$ws.Range("S15").Value = 0.006201851476855156
$ws.Range("T15").Value = 0.004184266679723272

$ws.Range("G16").Value = 473.968811
$ws.Range("H16").Value = 1421.906433
$ws.Range("I16").Value = 0.3582615042098434
$ws.Range("J16").Value = 0.360475806319893
$ws.Range("M16").Value = 104.9290516666667
$ws.Range("N16").Value = 314.787155
$ws.Range("O16").Value = 0.4249782927033401
$ws.Range("P16").Value = 0.4274447862635018
$ws.Range("Q16").Value = 49733.09785780757
$ws.Range("R16").Value = 447597.8807202681
$ws.Range("S16").Value = 0.1522533624004297
$ws.Range("T16").Value = 0.1540835039855701

$ws.Range("G17").Value = 24.3798835
$ws.Range("H17").Value = 48.759767
$ws.Range("I17").Value = 0.01842816137361988
$ws.Range("J17").Value = 0.01236137337687614
$ws.Range("M17").Value = 11.53121
$ws.Range("N17").Value = 34.59363
$ws.Range("O17").Value = 0.04670311854310271
$ws.Range("P17").Value = 0.04697417460197403
$ws.Range("Q17").Value = 281.129556414035
$ws.Range("R17").Value = 1686.77733848421
$ws.Range("S17").Value = 0.0008606526051635959
$ws.Range("T17").Value = 0.000580665311325573

$ws.Range("G18").Value = 24.3798835
$ws.Range("H18").Value = 48.759767
$ws.Range("I18").Value = 0.01842816137361988
$ws.Range("J18").Value = 0.01236137337687614
$ws.Range("O18").Value = 0.04941252241252041
$ws.Range("P18").Value = 0.04969930333854504
$ws.Range("Q18").Value = 297.438820799729
$ws.Range("R18").Value = 1784.632924798374
$ws.Range("S18").Value = 0.0009105819368955355
$ws.Range("T18").Value = 0.0006143516451383819

$ws.Range("G19").Value = 24.3798835
$ws.Range("H19").Value = 48.759767
$ws.Range("I19").Value = 0.01842816137361988
$ws.Range("J19").Value = 0.01236137337687614
$ws.Range("M19").Value = 113.969907
$ws.Range("N19").Value = 341.909721
$ws.Range("O19").Value = 0.4615951038067463
$ws.Range("P19").Value = 0.4642741144067919
$ws.Range("Q19").Value = 2778.573055165834
$ws.Range("R19").Value = 16671.43833099501
$ws.Range("S19").Value = 0.008506349062223542
$ws.Range("T19").Value = 0.005739065677400863

$ws.Range("G20").Value = 24.3798835
$ws.Range("H20").Value = 48.759767
$ws.Range("I20").Value = 0.01842816137361988
$ws.Range("J20").Value = 0.01236137337687614
$ws.Range("M20").Value = 4.2741545
$ws.Range("N20").Value = 8.548309
$ws.Range("O20").Value = 0.01731096253429049
$ws.Range("P20").Value = 0.01160762138918714
$ws.Range("Q20").Value = 104.2033887710008
$ws.Range("R20").Value = 416.813555084003
$ws.Range("S20").Value = 0.000319009211114593
$ws.Range("T20").Value = 0.000143486142009156

$ws.Range("G21").Value = 24.3798835
$ws.Range("H21").Value = 48.759767
$ws.Range("I21").Value = 0.01842816137361988
$ws.Range("J21").Value = 0.01236137337687614
$ws.Range("M21").Value = 104.9290516666667
$ws.Range("N21").Value = 314.787155
$ws.Range("O21").Value = 0.4249782927033401
$ws.Range("P21").Value = 0.4274447862635018
$ws.Range("Q21").Value = 2558.158055398814
$ws.Range("R21").Value = 15348.94833239289
$ws.Range("S21").Value = 0.007831568558222616
$ws.Range("T21").Value = 0.005283804601002162

$ws.Range("G22").Value = 692.3345543333334
$ws.Range("H22").Value = 2077.003663
$ws.Range("I22").Value = 0.5233188621179371
$ws.Range("J22").Value = 0.5265533320428379
$ws.Range("M22").Value = 11.53121
$ws.Range("N22").Value = 34.59363
$ws.Range("O22").Value = 0.04670311854310271
$ws.Range("P22").Value = 0.04697417460197403
$ws.Range("Q22").Value = 7983.455136274077
$ws.Range("R22").Value = 71851.09622646669
$ws.Range("S22").Value = 0.02444062285333564
$ws.Range("T22").Value = 0.02473440815663148

$ws.Range("G23").Value = 692.3345543333334
$ws.Range("H23").Value = 2077.003663
$ws.Range("I23").Value = 0.5233188621179371
$ws.Range("J23").Value = 0.5265533320428379
$ws.Range("O23").Value = 0.04941252241252041
$ws.Range("P23").Value = 0.04969930333854504
$ws.Range("Q23").Value = 8446.60202907912
$ws.Range("R23").Value = 76019.41826171208
$ws.Range("S23").Value = 0.02585850500329725
$ws.Range("T23").Value = 0.02616933377311863

$ws.Range("G24").Value = 692.3345543333334
$ws.Range("H24").Value = 2077.003663
$ws.Range("I24").Value = 0.5233188621179371
$ws.Range("J24").Value = 0.5265533320428379
$ws.Range("M24").Value = 113.969907
$ws.Range("N24").Value = 341.909721
$ws.Range("O24").Value = 0.4615951038067463
$ws.Range("P24").Value = 0.4642741144067919
$ws.Range("Q24").Value = 78905.30477025645
$ws.Range("R24").Value = 710147.742932308
$ws.Range("S24").Value = 0.2415614244833575
$ws.Range("T24").Value = 0.244465081922134

$ws.Range("G25").Value = 692.3345543333334
$ws.Range("H25").Value = 2077.003663
$ws.Range("I25").Value = 0.5233188621179371
$ws.Range("J25").Value = 0.5265533320428379
$ws.Range("M25").Value = 4.2741545
$ws.Range("N25").Value = 8.548309
$ws.Range("O25").Value = 0.01731096253429049
$ws.Range("P25").Value = 0.01160762138918714
$ws.Range("Q25").Value = 2959.144850909311
$ws.Range("R25").Value = 17754.86910545587
$ws.Range("S25").Value = 0.009059153215611142
$ws.Range("T25").Value = 0.006112031719568207

$ws.Range("G26").Value = 692.3345543333334
$ws.Range("H26").Value = 2077.003663
$ws.Range("I26").Value = 0.5233188621179371
$ws.Range("J26").Value = 0.5265533320428379
$ws.Range("M26").Value = 104.9290516666667
$ws.Range("N26").Value = 314.787155
$ws.Range("O26").Value = 0.4249782927033401
$ws.Range("P26").Value = 0.4274447862635018
$ws.Range("Q26").Value = 72646.00822226098
$ws.Range("R26").Value = 653814.0740003487
$ws.Range("S26").Value = 0.2223991565623356
$ws.Range("T26").Value = 0.2250724764713856
